# "Generate Report for Handback"
#
# The localization-status report previously only listed the outbound
# handoff info (source file + handoff xlf). This change fills in the
# handback columns once a translated file has come back "in sync" with
# en-US: the new "Latest Target File" (F) / "Latest Handback File" (G)
# columns get populated (as hyperlinked file names, matching the existing
# A/D columns), the Status text flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and the "Latest Handback DateTime"
# (H) is stamped with the real handback time instead of the zero-date
# placeholder.

$wb = $excel.ActiveWorkbook

$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d396d708eea43da2da5a41e2ce49398961c6f9ba/e2e/a.md"

$statusHandedBack = "Handed back: in sync with en-US"

# Cornflower blue (matches the workbook's existing HyperLink style), as a
# COM BGR color value.
$hyperlinkColor = 15570276

function Apply-HandbackRow {
    param(
        $ws,
        [int]$row,
        [string]$targetXlfName,
        [string]$targetXlfUrl,
        [string]$handbackDatetime
    )

    # Status -> handed back, in sync with en-US.
    $ws.Range("C$row").Value2 = $statusHandedBack

    # Latest Handback DateTime.
    $ws.Range("H$row").Value2 = $handbackDatetime

    # Latest Target File (F) - hyperlinked, mirrors the "a.md" source link.
    $fCell = $ws.Range("F$row")
    $fCell.Value2 = "a.md"
    $ws.Hyperlinks.Add($fCell, $aMdUrl, "", "", "a.md")
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # Latest Handback File (G) - hyperlinked, mirrors the handoff xlf link.
    $gCell = $ws.Range("G$row")
    $gCell.Value2 = $targetXlfName
    $ws.Hyperlinks.Add($gCell, $targetXlfUrl, "", "", $targetXlfName)
    $gCell.Font.Underline = 2
    $gCell.Font.Color = $hyperlinkColor
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhCnXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ecdd39d9006a65edfce2b3a71a83c87d458a568/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

Apply-HandbackRow $wsZh 2 $zhCnXlfName $zhCnXlfUrl "2016-03-21 14:33:00"
Apply-HandbackRow $wsZh 3 $zhCnXlfName $zhCnXlfUrl "2016-03-21 14:33:00"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$deDeXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d340a8b62c1ee27895b9cd1289adf0185e01350/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

Apply-HandbackRow $wsDe 2 $deDeXlfName $deDeXlfUrl "2016-03-21 14:33:09"
Apply-HandbackRow $wsDe 3 $deDeXlfName $deDeXlfUrl "2016-03-21 14:33:09"

# ---- Overview sheet ----
# Shares the same "Status" text via the shared-string table, so it flips
# to "Handed back: in sync with en-US" along with the per-language sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $statusHandedBack
$wsOverview.Range("C2").Value2 = $statusHandedBack
$wsOverview.Range("B3").Value2 = $statusHandedBack
$wsOverview.Range("C3").Value2 = $statusHandedBack

Write-Output "Handback report generated."
